# Add "Raw and Clean Data from SSA for June 14th" as a new row (row 15) to the
# historical log sheet, and normalize the date format of the previous last
# row (row 14, column B) to match the rest of the date column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Propagate existing cell formatting down to the new row (row 15) ---
# Copy B14's current number format (date-only, "YYYY-MM-DD") to B15 *before*
# B14's own format is changed below, since the new last row should carry the
# format the previous last row used to have.
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)  # xlPasteFormats

# Copy A14's style (bold / bordered / centered) to A15.
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)  # xlPasteFormats

# Copy C14:G14's (default/unstyled) formatting to C15:G15.
$ws.Range("C14:G14").Copy()
$ws.Range("C15:G15").PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = $false

# --- 2. Normalize B14's format to match B2:B13 (datetime format) ---
$ws.Range("B14").NumberFormat = $ws.Range("B13").NumberFormat

# --- 3. Populate the new row of data (June 14th) ---
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 43996
$ws.Range("C15").Value = 146837
$ws.Range("D15").Value = 207076
$ws.Range("E15").Value = 52636
$ws.Range("F15").Value = 17141
$ws.Range("G15").Value = 32.5
